$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells(28, 8).Value = 610.4737
$ws.Cells(28, 9).Value = 603.8
$ws.Cells(28, 10).Value = 617.8889
$ws.Cells(28, 11).Value = 603.8
$ws.Cells(28, 12).Value = 617.8889
$ws.Cells(28, 13).Value = -118.8
$ws.Cells(28, 14).Value = -1587.8889
# Row 43
$ws.Cells(43, 8).Value = 1131.3636
$ws.Cells(43, 10).Value = 944.5
$ws.Cells(43, 12).Value = 944.5
$ws.Cells(43, 14).Value = -1082.5
# Row 62
$ws.Cells(62, 8).Value = 4297.5
$ws.Cells(62, 9).Value = 4296.6665
$ws.Cells(62, 11).Value = 4296.6665
$ws.Cells(62, 13).Value = -3672.6665
# Row 65
$ws.Cells(65, 8).Value = 4297.5
$ws.Cells(65, 9).Value = 4296.6665
$ws.Cells(65, 11).Value = 21483.3325
$ws.Cells(65, 13).Value = -18363.3325
# Row 127
$ws.Cells(127, 8).Value = 1165.5
$ws.Cells(127, 9).Value = 620
$ws.Cells(127, 10).Value = 1468.5555
$ws.Cells(127, 11).Value = 1860
$ws.Cells(127, 12).Value = 4405.666499999999
$ws.Cells(127, 13).Value = 3100
$ws.Cells(127, 14).Value = -14325.6665
# Row 131
$ws.Cells(131, 8).Value = 5380
$ws.Cells(131, 9).Value = 760
$ws.Cells(131, 11).Value = 2280
$ws.Cells(131, 13).Value = 2760
# Row 137
$ws.Cells(137, 8).Value = 2970.5
$ws.Cells(137, 9).Value = 3197.9333
$ws.Cells(137, 10).Value = 1833.3334
$ws.Cells(137, 11).Value = 9593.7999
$ws.Cells(137, 12).Value = 5500.0002
$ws.Cells(137, 13).Value = -7043.7999
$ws.Cells(137, 14).Value = -10600.0002
# Row 138
$ws.Cells(138, 8).Value = 1589.3125
$ws.Cells(138, 9).Value = 1126
$ws.Cells(138, 10).Value = 2773.3333
$ws.Cells(138, 11).Value = 3378
$ws.Cells(138, 12).Value = 8319.999899999999
$ws.Cells(138, 13).Value = 1762
$ws.Cells(138, 14).Value = -18599.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells(32, 8).Value = 8996.489
$ws.Cells(32, 9).Value = 11088.571
$ws.Cells(32, 10).Value = 1674.2
$ws.Cells(32, 11).Value = 11088.571
$ws.Cells(32, 12).Value = 1674.2
$ws.Cells(32, 13).Value = -10801.571
$ws.Cells(32, 14).Value = -2248.2
# Row 45
$ws.Cells(45, 8).Value = 4627.077
$ws.Cells(45, 9).Value = 4299.1875
$ws.Cells(45, 10).Value = 5151.7
$ws.Cells(45, 11).Value = 4299.1875
$ws.Cells(45, 12).Value = 5151.7
$ws.Cells(45, 13).Value = -3922.1875
$ws.Cells(45, 14).Value = -5905.7
# Row 61
$ws.Cells(61, 8).Value = 2978.3704
$ws.Cells(61, 9).Value = 1886.6
$ws.Cells(61, 10).Value = 4343.0835
$ws.Cells(61, 11).Value = 1886.6
$ws.Cells(61, 12).Value = 4343.0835
$ws.Cells(61, 13).Value = -1674.6
$ws.Cells(61, 14).Value = -4767.0835
# Row 122
$ws.Cells(122, 8).Value = 2862.125
$ws.Cells(122, 9).Value = 1719.6
$ws.Cells(122, 10).Value = 20000
$ws.Cells(122, 11).Value = 5158.799999999999
$ws.Cells(122, 12).Value = 60000
$ws.Cells(122, 13).Value = -2708.799999999999
$ws.Cells(122, 14).Value = -64900
# Row 132
$ws.Cells(132, 8).Value = 4090.96
$ws.Cells(132, 9).Value = 1694.3572
$ws.Cells(132, 10).Value = 7141.1816
$ws.Cells(132, 11).Value = 5083.071599999999
$ws.Cells(132, 12).Value = 21423.5448
$ws.Cells(132, 13).Value = -2553.071599999999
$ws.Cells(132, 14).Value = -26483.5448
# Row 136
$ws.Cells(136, 8).Value = 2978.3704
$ws.Cells(136, 9).Value = 1886.6
$ws.Cells(136, 10).Value = 4343.0835
$ws.Cells(136, 11).Value = 5659.799999999999
$ws.Cells(136, 12).Value = 13029.2505
$ws.Cells(136, 13).Value = -3109.799999999999
$ws.Cells(136, 14).Value = -18129.2505

$ws = $wb.Worksheets.Item("BSM")
# Row 133
$ws.Cells(133, 8).Value = 0
$ws.Cells(133, 10).Value = 0
$ws.Cells(133, 12).Value = 0
$ws.Cells(133, 14).ClearContents()
# Row 134
$ws.Cells(134, 8).Value = 5849.231
$ws.Cells(134, 9).Value = 2848.7368
$ws.Cells(134, 10).Value = 8699.700000000001
$ws.Cells(134, 11).Value = 8546.2104
$ws.Cells(134, 12).Value = 26099.1
$ws.Cells(134, 13).Value = -6011.2104
$ws.Cells(134, 14).Value = -31169.1

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells(16, 8).Value = 3457.625
$ws.Cells(16, 9).Value = 4470.3335
$ws.Cells(16, 10).Value = 2850
$ws.Cells(16, 11).Value = 4470.3335
$ws.Cells(16, 12).Value = 2850
$ws.Cells(16, 13).Value = -4183.3335
$ws.Cells(16, 14).Value = -3424
# Row 31
$ws.Cells(31, 8).Value = 2658.9744
$ws.Cells(31, 9).Value = 2017.3914
$ws.Cells(31, 10).Value = 3581.25
$ws.Cells(31, 11).Value = 2017.3914
$ws.Cells(31, 12).Value = 3581.25
$ws.Cells(31, 13).Value = -1722.3914
$ws.Cells(31, 14).Value = -4171.25
# Row 34
$ws.Cells(34, 8).Value = 2658.9744
$ws.Cells(34, 9).Value = 2017.3914
$ws.Cells(34, 10).Value = 3581.25
$ws.Cells(34, 11).Value = 2017.3914
$ws.Cells(34, 12).Value = 3581.25
$ws.Cells(34, 13).Value = -1815.3914
$ws.Cells(34, 14).Value = -3985.25
# Row 107
$ws.Cells(107, 8).Value = 1202.9375
$ws.Cells(107, 9).Value = 558.05
$ws.Cells(107, 10).Value = 2277.75
$ws.Cells(107, 11).Value = 558.05
$ws.Cells(107, 12).Value = 2277.75
$ws.Cells(107, 13).Value = 1361.95
$ws.Cells(107, 14).Value = -6117.75
# Row 113
$ws.Cells(113, 8).Value = 3457.625
$ws.Cells(113, 9).Value = 4470.3335
$ws.Cells(113, 10).Value = 2850
$ws.Cells(113, 11).Value = 4470.3335
$ws.Cells(113, 12).Value = 2850
$ws.Cells(113, 13).Value = -2300.3335
$ws.Cells(113, 14).Value = -7190
# Row 132
$ws.Cells(132, 8).Value = 3641.8333
$ws.Cells(132, 9).Value = 1784
$ws.Cells(132, 10).Value = 5499.6665
$ws.Cells(132, 11).Value = 5352
$ws.Cells(132, 12).Value = 16498.9995
$ws.Cells(132, 13).Value = -2822
$ws.Cells(132, 14).Value = -21558.9995

$ws = $wb.Worksheets.Item("CUL")
# Row 101
$ws.Cells(101, 8).Value = 9333.333000000001
$ws.Cells(101, 10).Value = 9333.333000000001
$ws.Cells(101, 12).Value = 27999.999
$ws.Cells(101, 14).Value = -32867.999

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Cells(126, 8).Value = 14709330
$ws.Cells(126, 9).Value = 17860650
$ws.Cells(126, 10).Value = 3165.6667
$ws.Cells(126, 11).Value = 53581950
$ws.Cells(126, 12).Value = 9497.000100000001
$ws.Cells(126, 13).Value = -53579480
$ws.Cells(126, 14).Value = -14437.0001
# Row 132
$ws.Cells(132, 8).Value = 3265.818
$ws.Cells(132, 9).Value = 5956
$ws.Cells(132, 10).Value = 2668
$ws.Cells(132, 11).Value = 17868
$ws.Cells(132, 12).Value = 8004
$ws.Cells(132, 13).Value = -15338
$ws.Cells(132, 14).Value = -13064

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Cells(122, 8).Value = 10711.883
$ws.Cells(122, 9).Value = 17668.857
$ws.Cells(122, 10).Value = 5842
$ws.Cells(122, 11).Value = 53006.571
$ws.Cells(122, 12).Value = 17526
$ws.Cells(122, 13).Value = -50556.571
$ws.Cells(122, 14).Value = -22426
# Row 132
$ws.Cells(132, 8).Value = 71434600
$ws.Cells(132, 9).Value = 200011360
$ws.Cells(132, 10).Value = 3066.2222
$ws.Cells(132, 11).Value = 600034080
$ws.Cells(132, 12).Value = 9198.6666
$ws.Cells(132, 13).Value = -600031550
$ws.Cells(132, 14).Value = -14258.6666
# Row 136
$ws.Cells(136, 8).Value = 17858570
$ws.Cells(136, 9).Value = 23810880
$ws.Cells(136, 10).Value = 1641
$ws.Cells(136, 11).Value = 71432640
$ws.Cells(136, 12).Value = 4923
$ws.Cells(136, 13).Value = -71430090
$ws.Cells(136, 14).Value = -10023

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Cells(122, 8).Value = 9887.046
$ws.Cells(122, 9).Value = 13568.272
$ws.Cells(122, 10).Value = 6205.8184
$ws.Cells(122, 11).Value = 40704.81600000001
$ws.Cells(122, 12).Value = 18617.4552
$ws.Cells(122, 13).Value = -38254.81600000001
$ws.Cells(122, 14).Value = -23517.4552
# Row 126
$ws.Cells(126, 8).Value = 9713.532999999999
$ws.Cells(126, 9).Value = 10900.23
$ws.Cells(126, 10).Value = 2000
$ws.Cells(126, 11).Value = 32700.69
$ws.Cells(126, 12).Value = 6000
$ws.Cells(126, 13).Value = -30230.69
$ws.Cells(126, 14).Value = -10940
# Row 132
$ws.Cells(132, 8).Value = 6647
$ws.Cells(132, 9).Value = 7442
$ws.Cells(132, 10).Value = 6249.5
$ws.Cells(132, 11).Value = 22326
$ws.Cells(132, 12).Value = 18748.5
$ws.Cells(132, 13).Value = -19796
$ws.Cells(132, 14).Value = -23808.5
# Row 136
$ws.Cells(136, 8).Value = 13891010
$ws.Cells(136, 9).Value = 25001188
$ws.Cells(136, 10).Value = 3287.5
$ws.Cells(136, 11).Value = 75003564
$ws.Cells(136, 12).Value = 9862.5
$ws.Cells(136, 13).Value = -75001014
$ws.Cells(136, 14).Value = -14962.5

